$d = $word.ActiveDocument

# "... se prikazuju odgovarajući način ..." -> "... se prikazuju na odgovarajući način ..."
$d.Content.Find.Execute("prikazuju", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "prikazuju na", 2) | Out-Null
